$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = 0

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 436.75
$ws.Range("J19").Value = 478.5
$ws.Range("L19").Value = 478.5
$ws.Range("N19").Value = -828.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 8396.471
$ws.Range("I113").Value = 7286.091
$ws.Range("K113").Value = 7286.091
$ws.Range("M113").Value = -4032.091

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4998
$ws.Range("I2").Value = 5247.5835
$ws.Range("J2").Value = 4249.25
$ws.Range("K2").Value = 5247.5835
$ws.Range("L2").Value = 4249.25
$ws.Range("M2").Value = -5134.5835
$ws.Range("N2").Value = -4475.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7453.533
$ws.Range("J32").Value = 13333.333
$ws.Range("L32").Value = 13333.333
$ws.Range("N32").Value = -13907.333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1745.9487
$ws.Range("I74").Value = 689.1
$ws.Range("K74").Value = 689.1
$ws.Range("M74").Value = 184.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1745.9487
$ws.Range("I77").Value = 689.1
$ws.Range("K77").Value = 3445.5
$ws.Range("M77").Value = 922.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 4998
$ws.Range("I116").Value = 5247.5835
$ws.Range("J116").Value = 4249.25
$ws.Range("K116").Value = 5247.5835
$ws.Range("L116").Value = 4249.25
$ws.Range("M116").Value = -2953.5835
$ws.Range("N116").Value = -8837.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4998
$ws.Range("I3").Value = 5247.5835
$ws.Range("J3").Value = 4249.25
$ws.Range("K3").Value = 5247.5835
$ws.Range("L3").Value = 4249.25
$ws.Range("M3").Value = -5133.5835
$ws.Range("N3").Value = -4477.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6388.96
$ws.Range("I86").Value = 6870.524
$ws.Range("K86").Value = 6870.524
$ws.Range("M86").Value = -5747.524

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 6388.96
$ws.Range("I89").Value = 6870.524
$ws.Range("K89").Value = 34352.62
$ws.Range("M89").Value = -28736.62

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 21138.666
$ws.Range("I99").Value = 23340
$ws.Range("K99").Value = 23340
$ws.Range("M99").Value = -21842

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 6577.68
$ws.Range("I105").Value = 10008.833
$ws.Range("J105").Value = 3410.4614
$ws.Range("K105").Value = 10008.833
$ws.Range("L105").Value = 3410.4614
$ws.Range("M105").Value = -8261.833000000001
$ws.Range("N105").Value = -6904.4614

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1252
$ws.Range("I107").Value = 1252
$ws.Range("K107").Value = 1252
$ws.Range("M107").Value = 668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 171054740
$ws.Range("J108").Value = 171054740
$ws.Range("L108").Value = 171054740
$ws.Range("N108").Value = -171062420

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2527.6428
$ws.Range("I134").Value = 1322.5714
$ws.Range("K134").Value = 3967.7142
$ws.Range("M134").Value = -1432.7142

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4088
$ws.Range("J31").Value = 5329.7
$ws.Range("L31").Value = 5329.7
$ws.Range("N31").Value = -5919.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4088
$ws.Range("J34").Value = 5329.7
$ws.Range("L34").Value = 5329.7
$ws.Range("N34").Value = -5733.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 25374.75
$ws.Range("I39").Value = 29333
$ws.Range("J39").Value = 13500
$ws.Range("K39").Value = 29333
$ws.Range("L39").Value = 13500
$ws.Range("M39").Value = -28942
$ws.Range("N39").Value = -14282

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H49").Value = 25374.75
$ws.Range("I49").Value = 29333
$ws.Range("J49").Value = 13500
$ws.Range("K49").Value = 29333
$ws.Range("L49").Value = 13500
$ws.Range("M49").Value = -29151
$ws.Range("N49").Value = -13864

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2614.2703
$ws.Range("I58").Value = 2164.1333
$ws.Range("J58").Value = 4543.4287
$ws.Range("K58").Value = 2164.1333
$ws.Range("L58").Value = 4543.4287
$ws.Range("M58").Value = -1961.1333
$ws.Range("N58").Value = -4949.4287

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 69781.75
$ws.Range("J59").Value = 69781.75
$ws.Range("L59").Value = 69781.75
$ws.Range("N59").Value = -72071.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4747.7036
$ws.Range("I132").Value = 5007.8335
$ws.Range("K132").Value = 15023.5005
$ws.Range("M132").Value = -12493.5005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1469.45
$ws.Range("I134").Value = 993.17645
$ws.Range("J134").Value = 4168.3335
$ws.Range("K134").Value = 2979.52935
$ws.Range("L134").Value = 12505.0005
$ws.Range("M134").Value = -444.5293500000002
$ws.Range("N134").Value = -17575.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 110778
$ws.Range("J135").Value = 110778
$ws.Range("L135").Value = 110778
$ws.Range("N135").Value = -120918

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2614.2703
$ws.Range("I136").Value = 2164.1333
$ws.Range("J136").Value = 4543.4287
$ws.Range("K136").Value = 6492.3999
$ws.Range("L136").Value = 13630.2861
$ws.Range("M136").Value = -3942.3999
$ws.Range("N136").Value = -18730.2861

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 213926.53
$ws.Range("J141").Value = 224986.42
$ws.Range("L141").Value = 224986.42
$ws.Range("N141").Value = -235346.42

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 4810.143
$ws.Range("I18").Value = 6467.8
$ws.Range("K18").Value = 19403.4
$ws.Range("M18").Value = -19234.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 271999.4
$ws.Range("I128").Value = 271999.4
$ws.Range("K128").Value = 815998.2000000001
$ws.Range("M128").Value = -811018.2000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3650.0334
$ws.Range("I131").Value = 4227.048
$ws.Range("J131").Value = 2303.6667
$ws.Range("K131").Value = 12681.144
$ws.Range("L131").Value = 6911.000100000001
$ws.Range("M131").Value = -7641.144
$ws.Range("N131").Value = -16991.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 2224948.2
$ws.Range("I139").Value = 4445388.5
$ws.Range("J139").Value = 4507.778
$ws.Range("K139").Value = 13336165.5
$ws.Range("L139").Value = 13523.334
$ws.Range("M139").Value = -13331025.5
$ws.Range("N139").Value = -23803.334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1003909.9
$ws.Range("I140").Value = 1003909.9
$ws.Range("K140").Value = 3011729.7
$ws.Range("M140").Value = -3006549.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 17112.5
$ws.Range("I122").Value = 14350.75
$ws.Range("K122").Value = 43052.25
$ws.Range("M122").Value = -40602.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 17817.072
$ws.Range("J126").Value = 16913.166
$ws.Range("L126").Value = 50739.49800000001
$ws.Range("N126").Value = -55679.49800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2945.2942
$ws.Range("J132").Value = 2571.4285
$ws.Range("L132").Value = 7714.2855
$ws.Range("N132").Value = -12774.2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2259.5356
$ws.Range("I16").Value = 2118.6956
$ws.Range("J16").Value = 2907.4
$ws.Range("K16").Value = 2118.6956
$ws.Range("L16").Value = 2907.4
$ws.Range("M16").Value = -1948.6956
$ws.Range("N16").Value = -3247.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1517.875
$ws.Range("I46").Value = 777.44446
$ws.Range("K46").Value = 777.44446
$ws.Range("M46").Value = -589.44446

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H92").Value = 58999
$ws.Range("J92").Value = 58999
$ws.Range("L92").Value = 58999
$ws.Range("N92").Value = -63991

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3164.9
$ws.Range("I93").Value = 4042.7896
$ws.Range("K93").Value = 4042.7896
$ws.Range("M93").Value = -2794.7896

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6134.433
$ws.Range("I122").Value = 6235.4375
$ws.Range("J122").Value = 6019
$ws.Range("K122").Value = 18706.3125
$ws.Range("L122").Value = 18057
$ws.Range("M122").Value = -16256.3125
$ws.Range("N122").Value = -22957

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 93997
$ws.Range("J46").Value = 93997
$ws.Range("L46").Value = 93997
$ws.Range("N46").Value = -94459

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 7119.068
$ws.Range("I132").Value = 8404.091
$ws.Range("K132").Value = 25212.273
$ws.Range("M132").Value = -22682.273

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H134").Value = 93997
$ws.Range("J134").Value = 93997
$ws.Range("L134").Value = 281991
$ws.Range("N134").Value = -287061

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 387040.97
$ws.Range("I136").Value = 418077.12
$ws.Range("K136").Value = 1254231.36
$ws.Range("M136").Value = -1251681.36
